# Plantilla Visitas_TEK.xlsx — daily visitor log update
#
# Changes applied (per commit "update index page layout ... add download
# template feature ..." as reflected in this workbook's data):
#   1. RUT header note (D3) is simplified from a rich-text run to a single
#      plain-text string and the quoting in the instructions changes from
#      sin "." to sin "".
#   2. The "Fecha de ingreso" date for the first four visitor rows moves
#      from 12/19/2025 to 12/22/2025 (next business day).
#   3. Row 8's visitor is replaced: Patricio Sanchez (RUT 15727383-3) is
#      swapped out for Hector Tallaedo (RUT 22309814-2); the vehicle
#      plate / motive / date columns for that row are left blank since the
#      new visitor doesn't drive in, and that row is restyled to match the
#      normal data rows instead of the old "extra visitor" styling.
#   4. Rows 9 and 10 (previously Osvaldo Carrasco and Gonzalo Quezada) are
#      cleared out entirely — no more extra visitors/vehicles.
#   5. Selection cursor left on I7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. RUT column header: collapse the rich-text note into plain text and
#    drop the period inside the quotes ("." -> "").
$ws.Range("D3").Value = "RUT`n(Ingresar sin """" y con digito verificador)"

# 2. Push the ingreso date for the first four rows to 12/22/2025.
$ws.Range("I4:I7").Value = 46013

# 3. Row 8 becomes Hector Tallaedo, with no vehicle / motive / date.
#    First, copy the normal-data-row formatting (from row 4) onto the
#    name/company cells so row 8 no longer looks like the old "extra
#    visitor" block.
$ws.Range("B4").Copy()
$ws.Range("B8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B4").Copy()
$ws.Range("D8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D8").VerticalAlignment = -4108   # xlVAlignCenter

$ws.Range("B8").Value = "Héctor "
$ws.Range("C8").Value = "Tallaedo"
$ws.Range("D8").Value = "22309814-2"
$ws.Range("G8:I8").ClearContents()

# 4. Rows 9 and 10 (the other two extra visitors) are removed entirely.
$ws.Range("B9:I10").ClearContents()

# 5. Leave the selection on I7.
$ws.Range("I7").Select()
